$p = $ppt.ActivePresentation

# --- Update the cached "datetimeFigureOut" date placeholder text wherever
# it appears (the slide master and every slide layout) from 12/1/2020 to
# 12/6/2020, as PowerPoint does when it re-caches the auto date field.
$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        if ($shp.TextFrame.TextRange.Text -eq "12/1/2020") {
            $shp.TextFrame.TextRange.Text = "12/6/2020"
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq "12/1/2020") {
                $shp.TextFrame.TextRange.Text = "12/6/2020"
            }
        }
    }
}

# --- Slide 9: rename the "Extending the Botnet" section-header rectangle to
# "Removing the Infection" and widen it to fit the new (longer) text.
$s = $p.Slides.Item(9)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Extending the Botnet") {
        $shp.TextFrame.TextRange.Text = "Removing the Infection"
        $shp.Width = 3115159 / 12700
    }
}
